$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (headers) is untouched - its content and style ("s=1", bold+border)
# stay exactly as they are in the source workbook.

# New TPM-based data table (5 rows instead of 8), keyed by column letter so we
# can rewrite the sheet column-by-column (this mirrors how the source
# pipeline serialised the shared-string table: categorical column A first,
# sorted, then B, C, D, then the numeric metrics).
$colA = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")
$colB = @("Ccl12", "Ccl12", "Ccl12", "Ccl12", "Ccl12")
$colC = @("Ccr10", "Ccr10", "Ccr10", "Ccr10", "Ccr10")
$colD = @("MuSCs", "MuSCs", "MuSCs", "MuSCs", "MuSCs")
$colE = @(1, 1, 3, 1, 3)
$colF = @(0.3333333333333333, 0.3333333333333333, 1, 0.5, 1)
$colG = @(0.2401406666666667, 0.643692, 60.66100033333333, 1.6570225, 46.31492933333334)
$colH = @(0.720422, 1.931076, 181.983001, 3.314045, 138.944788)
$colI = @(0.002192729333974893, 0.005877564804149375, 0.553896833491318, 0.0151303063043872, 0.4229025660661706)
$colJ = @(0.002203844280310985, 0.005907358183739276, 0.5567045368793268, 0.01013800122420362, 0.4250462594324195)
$colK = @(1, 1, 1, 1, 1)
$colL = @(0.5, 0.5, 0.5, 0.5, 0.5)
$colM = @(0.1452775, 0.1452775, 0.1452775, 0.1452775, 0.1452775)
$colN = @(0.290555, 0.290555, 0.290555, 0.290555, 0.290555)
$colO = @(1, 1, 1, 1, 1)
$colP = @(1, 1, 1, 1, 1)
$colQ = @(0.03488703570166667, 0.09351396453000001, 8.812678475925834, 0.24072808624375, 6.728517146223334)
$colR = @(0.20932221421, 0.5610837871800001, 52.876070855555, 0.9629123449750001, 40.37110287734001)
$colS = @(0.002192729333974893, 0.005877564804149375, 0.553896833491318, 0.0151303063043872, 0.4229025660661706)
$colT = @(0.002203844280310985, 0.005907358183739276, 0.5567045368793268, 0.01013800122420362, 0.4250462594324195)

$columns = @($colA, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI, $colJ, $colK, $colL, $colM, $colN, $colO, $colP, $colQ, $colR, $colS, $colT)

# Clear only the old data rows (2-9); row 1 (headers + style) is left intact.
$ws.Range("A2:T9").Clear()

for ($j = 0; $j -lt $columns.Count; $j++) {
    $col = $j + 1
    $colValues = $columns[$j]
    for ($i = 0; $i -lt $colValues.Count; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, $col).Value = $colValues[$i]
    }
}
